$wb = $excel.ActiveWorkbook

# --- ALC (Worksheets.Item(1)) ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(12,8).Value = 447.5
$ws.Cells.Item(12,9).Value = 293.75
$ws.Cells.Item(12,10).Value = 550
$ws.Cells.Item(12,11).Value = 293.75
$ws.Cells.Item(12,12).Value = 550
$ws.Cells.Item(12,13).Value = -123.75
$ws.Cells.Item(12,14).Value = -890
$ws.Cells.Item(28,8).Value = 834.8570999999999
$ws.Cells.Item(28,9).Value = 507.5
$ws.Cells.Item(28,10).Value = 1271.3334
$ws.Cells.Item(28,11).Value = 507.5
$ws.Cells.Item(28,12).Value = 1271.3334
$ws.Cells.Item(28,13).Value = -22.5
$ws.Cells.Item(28,14).Value = -2241.3334
$ws.Cells.Item(62,8).Value = 4407.457
$ws.Cells.Item(62,10).Value = 6170.8184
$ws.Cells.Item(62,12).Value = 6170.8184
$ws.Cells.Item(62,14).Value = -7418.8184
$ws.Cells.Item(65,8).Value = 4407.457
$ws.Cells.Item(65,10).Value = 6170.8184
$ws.Cells.Item(65,12).Value = 30854.092
$ws.Cells.Item(65,14).Value = -37094.092
$ws.Cells.Item(134,8).Value = 75485.164
$ws.Cells.Item(134,10).Value = 75485.164
$ws.Cells.Item(134,12).Value = 75485.164
$ws.Cells.Item(134,14).Value = -85625.164
$ws.Cells.Item(135,8).Value = 1338.4706
$ws.Cells.Item(135,9).Value = 1321.5834
$ws.Cells.Item(135,10).Value = 1379
$ws.Cells.Item(135,11).Value = 11894.2506
$ws.Cells.Item(135,12).Value = 12411
$ws.Cells.Item(135,13).Value = -9359.250599999999
$ws.Cells.Item(135,14).Value = -17481
$ws.Cells.Item(137,8).Value = 4404.4814
$ws.Cells.Item(137,9).Value = 1341.7778
$ws.Cells.Item(137,10).Value = 10529.889
$ws.Cells.Item(137,11).Value = 4025.3334
$ws.Cells.Item(137,12).Value = 31589.667
$ws.Cells.Item(137,13).Value = -1475.3334
$ws.Cells.Item(137,14).Value = -36689.667
$ws.Cells.Item(138,8).Value = 2744.1619
$ws.Cells.Item(138,9).Value = 1423
$ws.Cells.Item(138,11).Value = 4269
$ws.Cells.Item(138,13).Value = 871

# --- ARM (Worksheets.Item(2)) ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(24,8).Value = 111000
$ws.Cells.Item(24,10).Value = 111000
$ws.Cells.Item(24,12).Value = 111000
$ws.Cells.Item(24,14).Value = -111748
$ws.Cells.Item(32,8).Value = 7943074
$ws.Cells.Item(32,9).Value = 10419483
$ws.Cells.Item(32,11).Value = 10419483
$ws.Cells.Item(32,13).Value = -10419196
$ws.Cells.Item(45,8).Value = 38463800
$ws.Cells.Item(45,10).Value = 2538.5
$ws.Cells.Item(45,12).Value = 2538.5
$ws.Cells.Item(45,14).Value = -3292.5
$ws.Cells.Item(61,8).Value = 35720184
$ws.Cells.Item(61,9).Value = 31254414
$ws.Cells.Item(61,11).Value = 31254414
$ws.Cells.Item(61,13).Value = -31254202
$ws.Cells.Item(74,8).Value = 10009585
$ws.Cells.Item(74,9).Value = 19233640
$ws.Cells.Item(74,10).Value = 16858.916
$ws.Cells.Item(74,11).Value = 19233640
$ws.Cells.Item(74,12).Value = 16858.916
$ws.Cells.Item(74,13).Value = -19232766
$ws.Cells.Item(74,14).Value = -18606.916
$ws.Cells.Item(77,8).Value = 10009585
$ws.Cells.Item(77,9).Value = 19233640
$ws.Cells.Item(77,10).Value = 16858.916
$ws.Cells.Item(77,11).Value = 96168200
$ws.Cells.Item(77,12).Value = 84294.58
$ws.Cells.Item(77,13).Value = -96163832
$ws.Cells.Item(77,14).Value = -93030.58
$ws.Cells.Item(100,8).Value = 111000
$ws.Cells.Item(100,10).Value = 111000
$ws.Cells.Item(100,12).Value = 111000
$ws.Cells.Item(100,14).Value = -113164
$ws.Cells.Item(119,8).Value = 106499
$ws.Cells.Item(119,10).Value = 106499
$ws.Cells.Item(119,12).Value = 106499
$ws.Cells.Item(119,14).Value = -116175
$ws.Cells.Item(121,8).Value = 95127.5
$ws.Cells.Item(121,10).Value = 95127.5
$ws.Cells.Item(121,12).Value = 95127.5
$ws.Cells.Item(121,14).Value = -98621.5
$ws.Cells.Item(132,8).Value = 30967.756
$ws.Cells.Item(132,9).Value = 32398.213
$ws.Cells.Item(132,11).Value = 97194.639
$ws.Cells.Item(132,13).Value = -94664.639
$ws.Cells.Item(136,8).Value = 35720184
$ws.Cells.Item(136,9).Value = 31254414
$ws.Cells.Item(136,11).Value = 93763242
$ws.Cells.Item(136,13).Value = -93760692

# --- BSM (Worksheets.Item(3)) ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(134,8).Value = 241509.27
$ws.Cells.Item(134,9).Value = 1328.3513
$ws.Cells.Item(134,11).Value = 3985.0539
$ws.Cells.Item(134,13).Value = -1450.0539

# --- CRP (Worksheets.Item(4)) ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(7,8).Value = 253.41176
$ws.Cells.Item(7,9).Value = 215.54546
$ws.Cells.Item(7,10).Value = 322.83334
$ws.Cells.Item(7,11).Value = 215.54546
$ws.Cells.Item(7,12).Value = 322.83334
$ws.Cells.Item(7,13).Value = -102.54546
$ws.Cells.Item(7,14).Value = -548.83334
$ws.Cells.Item(70,8).Value = 0
$ws.Cells.Item(70,10).Value = 0
$ws.Cells.Item(70,12).Value = 0
$ws.Cells.Item(73,8).Value = 0
$ws.Cells.Item(73,10).Value = 0
$ws.Cells.Item(73,12).Value = 0
$ws.Cells.Item(122,8).Value = 1724.95
$ws.Cells.Item(122,9).Value = 1804.25
$ws.Cells.Item(122,11).Value = 5412.75
$ws.Cells.Item(122,13).Value = -2962.75
$ws.Cells.Item(133,8).Value = 55000
$ws.Cells.Item(133,10).Value = 55000
$ws.Cells.Item(133,12).Value = 55000
$ws.Cells.Item(133,14).Value = -60060
$ws.Cells.Item(134,8).Value = 3998.6
$ws.Cells.Item(134,9).Value = 1696.9
$ws.Cells.Item(134,11).Value = 5090.700000000001
$ws.Cells.Item(134,13).Value = -2555.700000000001
$ws.Cells.Item(141,8).Value = 349374.7
$ws.Cells.Item(141,10).Value = 377194.12
$ws.Cells.Item(141,12).Value = 377194.12
$ws.Cells.Item(141,14).Value = -387554.12
$ws.Cells.Item(70,14).ClearContents()
$ws.Cells.Item(73,14).ClearContents()

# --- CUL (Worksheets.Item(5)) ---
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(132,8).Value = 1916.1666
$ws.Cells.Item(132,9).Value = 1949.75
$ws.Cells.Item(132,11).Value = 17547.75
$ws.Cells.Item(132,13).Value = -15017.75
$ws.Cells.Item(134,8).Value = 10609.571
$ws.Cells.Item(134,10).Value = 12988.464
$ws.Cells.Item(134,12).Value = 38965.392
$ws.Cells.Item(134,14).Value = -49105.392
$ws.Cells.Item(136,8).Value = 2016
$ws.Cells.Item(136,9).Value = 2016
$ws.Cells.Item(136,11).Value = 6048
$ws.Cells.Item(136,13).Value = -948
$ws.Cells.Item(137,8).Value = 5880.364
$ws.Cells.Item(137,9).Value = 4733.3335
$ws.Cells.Item(137,10).Value = 6310.5
$ws.Cells.Item(137,11).Value = 14200.0005
$ws.Cells.Item(137,12).Value = 18931.5
$ws.Cells.Item(137,13).Value = -9100.000499999998
$ws.Cells.Item(137,14).Value = -29131.5
$ws.Cells.Item(141,8).Value = 260304.83
$ws.Cells.Item(141,9).Value = 505609.66
$ws.Cells.Item(141,11).Value = 1516828.98
$ws.Cells.Item(141,13).Value = -1511648.98

# --- GSM (Worksheets.Item(6)) ---
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(46,8).Value = 50000
$ws.Cells.Item(46,10).Value = 50000
$ws.Cells.Item(46,12).Value = 50000
$ws.Cells.Item(46,14).Value = -50312
$ws.Cells.Item(112,8).Value = 105875
$ws.Cells.Item(112,10).Value = 105875
$ws.Cells.Item(112,12).Value = 105875
$ws.Cells.Item(112,14).Value = -108091
$ws.Cells.Item(113,8).Value = 3958.238
$ws.Cells.Item(113,10).Value = 4077.25
$ws.Cells.Item(113,12).Value = 4077.25
$ws.Cells.Item(113,14).Value = -8417.25
$ws.Cells.Item(121,8).Value = 31500
$ws.Cells.Item(121,10).Value = 31500
$ws.Cells.Item(121,12).Value = 31500
$ws.Cells.Item(121,14).Value = -34994
$ws.Cells.Item(132,8).Value = 58826820
$ws.Cells.Item(132,9).Value = 71431920
$ws.Cells.Item(132,11).Value = 214295760
$ws.Cells.Item(132,13).Value = -214293230
$ws.Cells.Item(139,8).Value = 0
$ws.Cells.Item(139,10).Value = 0
$ws.Cells.Item(139,12).Value = 0
$ws.Cells.Item(139,14).ClearContents()

# --- LTW (Worksheets.Item(7)) ---
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(7,8).Value = 14496.211
$ws.Cells.Item(7,10).Value = 16013.454
$ws.Cells.Item(7,12).Value = 16013.454
$ws.Cells.Item(7,14).Value = -16237.454
$ws.Cells.Item(16,8).Value = 1003
$ws.Cells.Item(16,9).Value = 870.13336
$ws.Cells.Item(16,11).Value = 870.13336
$ws.Cells.Item(16,13).Value = -700.13336
$ws.Cells.Item(46,8).Value = 2496.2307
$ws.Cells.Item(46,9).Value = 1705.6666
$ws.Cells.Item(46,11).Value = 1705.6666
$ws.Cells.Item(46,13).Value = -1517.6666
$ws.Cells.Item(61,8).Value = 1274.3182
$ws.Cells.Item(61,9).Value = 983.5333000000001
$ws.Cells.Item(61,11).Value = 983.5333000000001
$ws.Cells.Item(61,13).Value = -781.5333000000001
$ws.Cells.Item(113,8).Value = 1274.3182
$ws.Cells.Item(113,9).Value = 983.5333000000001
$ws.Cells.Item(113,11).Value = 983.5333000000001
$ws.Cells.Item(113,13).Value = 1186.4667
$ws.Cells.Item(119,8).Value = 106000
$ws.Cells.Item(119,10).Value = 106000
$ws.Cells.Item(119,12).Value = 106000
$ws.Cells.Item(119,14).Value = -115676
$ws.Cells.Item(122,8).Value = 6285.9033
$ws.Cells.Item(122,9).Value = 5536.875
$ws.Cells.Item(122,10).Value = 7084.8667
$ws.Cells.Item(122,11).Value = 16610.625
$ws.Cells.Item(122,12).Value = 21254.6001
$ws.Cells.Item(122,13).Value = -14160.625
$ws.Cells.Item(122,14).Value = -26154.6001
$ws.Cells.Item(126,8).Value = 14496.211
$ws.Cells.Item(126,10).Value = 16013.454
$ws.Cells.Item(126,12).Value = 48040.362
$ws.Cells.Item(126,14).Value = -52980.362
$ws.Cells.Item(132,8).Value = 309049.97
$ws.Cells.Item(132,9).Value = 6454.8335
$ws.Cells.Item(132,10).Value = 3335001.2
$ws.Cells.Item(132,11).Value = 19364.5005
$ws.Cells.Item(132,12).Value = 10005003.6
$ws.Cells.Item(132,13).Value = -16834.5005
$ws.Cells.Item(132,14).Value = -10010063.6
$ws.Cells.Item(136,8).Value = 59232.477
$ws.Cells.Item(136,9).Value = 7785.5625
$ws.Cells.Item(136,11).Value = 23356.6875
$ws.Cells.Item(136,13).Value = -20806.6875

# --- WVR (Worksheets.Item(8)) ---
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(124,8).Value = 165656.67
$ws.Cells.Item(124,10).Value = 165656.67
$ws.Cells.Item(124,12).Value = 165656.67
$ws.Cells.Item(124,14).Value = -175476.67
$ws.Cells.Item(132,8).Value = 253527
$ws.Cells.Item(132,10).Value = 1673697.9
$ws.Cells.Item(132,12).Value = 5021093.699999999
$ws.Cells.Item(132,14).Value = -5026153.699999999
